{"js": "// Update the worksheet date and the 25 two-digit multiplication problems.\n// The document body contains exactly one title paragraph (the date) followed\n// by a single 20-row x 5-column table; every 5th row holds the problem text\n// (\"NN\u00d7NN=\"), the rows in between are blank spacer rows. `body.paragraphs`\n// walks the whole story (title + every table-cell paragraph) in document\n// order, so we can address each target paragraph positionally instead of\n// searching by text (a couple of the new values coincide with OLD values\n// elsewhere in the sheet, so positional addressing avoids any ambiguity).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// [paragraph index in body.paragraphs, new text]\nconst replacements = [\n  [0, \"2025-11-18 Tuesday\"],\n  [1, \"57\u00d739=\"],\n  [2, \"36\u00d748=\"],\n  [3, \"13\u00d765=\"],\n  [4, \"87\u00d753=\"],\n  [5, \"90\u00d717=\"],\n  [21, \"81\u00d789=\"],\n  [22, \"42\u00d769=\"],\n  [23, \"71\u00d797=\"],\n  [24, \"37\u00d723=\"],\n  [25, \"93\u00d784=\"],\n  [46, \"45\u00d768=\"],\n  [47, \"20\u00d764=\"],\n  [48, \"35\u00d722=\"],\n  [49, \"88\u00d766=\"],\n  [50, \"96\u00d732=\"],\n  [71, \"15\u00d779=\"],\n  [72, \"30\u00d779=\"],\n  [73, \"58\u00d790=\"],\n  [74, \"54\u00d744=\"],\n  [75, \"55\u00d757=\"],\n  [96, \"25\u00d779=\"],\n  [97, \"43\u00d712=\"],\n  [98, \"39\u00d716=\"],\n  [99, \"28\u00d786=\"],\n  [100, \"43\u00d774=\"],\n];\n\nfor (const [idx, newText] of replacements) {\n  paragraphs.items[idx].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems.\n# Each \"NN\u00d7NNN=\" prompt (and the date heading) is unique text in the\n# document, so a straightforward Find/Replace (one hit at a time, in\n# document order) locates every target unambiguously. Processing the\n# pairs in document order also keeps the single value that is both an\n# \"old\" and a \"new\" string in this edit (\"25\u00d779=\") unambiguous: its\n# original occurrence is consumed (replaced) before the later problem\n# is rewritten to reuse that same text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-11-17 Monday\"; New = \"2025-11-18 Tuesday\" },\n    @{ Old = \"95\u00d782=\"; New = \"57\u00d739=\" },\n    @{ Old = \"47\u00d749=\"; New = \"36\u00d748=\" },\n    @{ Old = \"67\u00d728=\"; New = \"13\u00d765=\" },\n    @{ Old = \"25\u00d779=\"; New = \"87\u00d753=\" },\n    @{ Old = \"26\u00d780=\"; New = \"90\u00d717=\" },\n    @{ Old = \"50\u00d788=\"; New = \"81\u00d789=\" },\n    @{ Old = \"89\u00d745=\"; New = \"42\u00d769=\" },\n    @{ Old = \"52\u00d754=\"; New = \"71\u00d797=\" },\n    @{ Old = \"35\u00d781=\"; New = \"37\u00d723=\" },\n    @{ Old = \"68\u00d765=\"; New = \"93\u00d784=\" },\n    @{ Old = \"22\u00d781=\"; New = \"45\u00d768=\" },\n    @{ Old = \"18\u00d775=\"; New = \"20\u00d764=\" },\n    @{ Old = \"39\u00d751=\"; New = \"35\u00d722=\" },\n    @{ Old = \"27\u00d754=\"; New = \"88\u00d766=\" },\n    @{ Old = \"27\u00d787=\"; New = \"96\u00d732=\" },\n    @{ Old = \"74\u00d769=\"; New = \"15\u00d779=\" },\n    @{ Old = \"72\u00d758=\"; New = \"30\u00d779=\" },\n    @{ Old = \"32\u00d750=\"; New = \"58\u00d790=\" },\n    @{ Old = \"53\u00d725=\"; New = \"54\u00d744=\" },\n    @{ Old = \"41\u00d747=\"; New = \"55\u00d757=\" },\n    @{ Old = \"29\u00d746=\"; New = \"25\u00d779=\" },\n    @{ Old = \"22\u00d789=\"; New = \"43\u00d712=\" },\n    @{ Old = \"42\u00d783=\"; New = \"39\u00d716=\" },\n    @{ Old = \"55\u00d787=\"; New = \"28\u00d786=\" },\n    @{ Old = \"63\u00d788=\"; New = \"43\u00d774=\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 1)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $($r.Old)\"\n    }\n}\n"}
